# Update view-count figures (column F) on two worksheets to match the
# latest scrape snapshot ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7117
$ws1.Range("F7").Value = 6436
$ws1.Range("F9").Value = 1997
$ws1.Range("F17").Value = 8245
$ws1.Range("F28").Value = 180
$ws1.Range("F32").Value = 433
$ws1.Range("F38").Value = 3944

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 7121
$ws4.Range("F11").Value = 6436
$ws4.Range("F13").Value = 1997
$ws4.Range("F23").Value = 8245
$ws4.Range("F32").Value = 180
$ws4.Range("F36").Value = 433
$ws4.Range("F44").Value = 3944
